# Đổi tên các sheet cho dễ theo dõi
$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "Version Control"
$wb.Worksheets.Item(2).Name = "Project Plan"
$wb.Worksheets.Item(3).Name = "Iteration 1"

$ws1 = $wb.Worksheets.Item("Version Control")
$ws2 = $wb.Worksheets.Item("Project Plan")
$ws3 = $wb.Worksheets.Item("Iteration 1")

# Update selection/cursor position on each sheet
$ws1.Range("C12").Select()
$ws3.Range("C12").Select()

# "Project Plan" ends up as the active/selected sheet (activeTab index 1)
$ws2.Activate()
$ws2.Range("C9").Select()
